# Removed Test Case Inter-Dependency
$wb = $excel.ActiveWorkbook

$wsInput  = $wb.Worksheets.Item("ProductLoanInput")
$wsOutput = $wb.Worksheets.Item("ProductLoanOutput")

# The product name used to be shared verbatim with another test case,
# causing the two automated runs to stomp on each other's product.
# Give this sheet's run its own distinct product name / short code so the
# two test cases no longer depend on each other's execution order.
$newProductName = "2595-RBI-EI-DB-DL-REC-NOCOM-RNI-CTPD-DL-MD-TR-2-DATE-VAR-INST-DISBURSE-FEE-%INT-1st"

$wsInput.Range("B1").Value = $newProductName
$wsOutput.Range("B1").Value = $newProductName

# shortname becomes an independent alphanumeric code instead of reusing the
# numeric product id that the other test also used.
$wsInput.Range("B2").Value = "259e"

# Collapse the old B2:B3 selection down to the single cell B2.
$wsInput.Range("B2").Select()

# Make the input sheet the active tab instead of the output sheet.
$wsInput.Activate()
